$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

# New values for rows 2-6 (and mirrored identically into rows 8-12)
$data = @{
    2  = @{ D = 0.0001685949973762035; E = 0.04236340383067727; G = 0.003296305891126394; H = 0.006117138545960188; I = 0.01170196989551187; J = 0.01729478221386671; K = 0.001073705032467842 }
    3  = @{ D = 0.002599451690912247;  E = 0.05051677720621228; G = 0.003446005284786224; H = 0.009040128905326128; I = 0.01246348908171058; J = 0.02132421731948853; K = 0.001154523342847824 }
    4  = @{ D = 0.004177008755505085;  E = 0.07502277009189129; G = 0.005170559510588646; H = 0.0129931503906846;   I = 0.01895631477236748; J = 0.03162738122045994; K = 0.001702926121652126 }
    5  = @{ D = 0.0003022640012204647; E = 0.0695937117561698;  G = 0.005592648405581713; H = 0.009667002130299807; I = 0.01953458599746227; J = 0.02783872047439218; K = 0.001880043651908636 }
    6  = @{ D = 0.004715193994343281;  E = 0.330788949970156;   G = 0.0101100062020123;   H = 0.02337119122967124;  I = 0.2516961521469057;  J = 0.03204407915472984; K = 0.003974684048444033 }
    8  = @{ D = 0.0001685949973762035; E = 0.04236340383067727; G = 0.003296305891126394; H = 0.006117138545960188; I = 0.01170196989551187; J = 0.01729478221386671; K = 0.001073705032467842 }
    9  = @{ D = 0.002599451690912247;  E = 0.05051677720621228; G = 0.003446005284786224; H = 0.009040128905326128; I = 0.01246348908171058; J = 0.02132421731948853; K = 0.001154523342847824 }
    10 = @{ D = 0.004177008755505085;  E = 0.07502277009189129; G = 0.005170559510588646; H = 0.0129931503906846;   I = 0.01895631477236748; J = 0.03162738122045994; K = 0.001702926121652126 }
    11 = @{ D = 0.0003022640012204647; E = 0.0695937117561698;  G = 0.005592648405581713; H = 0.009667002130299807; I = 0.01953458599746227; J = 0.02783872047439218; K = 0.001880043651908636 }
    12 = @{ D = 0.004715193994343281;  E = 0.330788949970156;   G = 0.0101100062020123;   H = 0.02337119122967124;  I = 0.2516961521469057;  J = 0.03204407915472984; K = 0.003974684048444033 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
